$wb = $excel.ActiveWorkbook

# Sheet 1 = 土地 (Land). Column I = property_category.
# Rows 2 and 3 had value "land"; change to "building".
$wsLand = $wb.Worksheets.Item(1)
$wsLand.Cells.Item(2, 9).Value = "building"
$wsLand.Cells.Item(3, 9).Value = "building"

# Sheet 3 = 汽車 (Car). Column H = property_category.
# Rows 2 and 3 had value "land"; change to "car".
$wsCar = $wb.Worksheets.Item(3)
$wsCar.Cells.Item(2, 8).Value = "car"
$wsCar.Cells.Item(3, 8).Value = "car"
